# feat: add 2022-Q4 data
#
# 1. The old "2022-Q3" fund-holdings sheet becomes "2022-Q4" (it keeps its
#    sheetId / rId, but its content is replaced by the new Q4 numbers).
# 2. A brand-new sheet named "2022-Q3" is inserted right after it, holding
#    the fund-holdings data that used to live in the "2022-Q3" sheet.
# 3. The "总计" (totals) sheet gets a new row for 2022-Q4, and the existing
#    2022-Q3 total row shifts down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: rename the existing "2022-Q3" sheet to "2022-Q4" and remember
# its current (soon to be old) fund-holdings content so we can recreate
# it verbatim on a brand-new "2022-Q3" sheet.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2022-Q3")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# Step 2: insert a new sheet right after it and copy the old Q3 data
# across before we overwrite the original sheet with the Q4 numbers.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $q4)
$q3.Name = "2022-Q3"
$q4.Range("A1:H6").Copy($q3.Range("A1"))

# ---------------------------------------------------------------------
# Step 3: overwrite the (renamed) Q4 sheet with the new fund data.
# ---------------------------------------------------------------------
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "000041"
$q4.Range("C2").Value = "华夏全球精选股票（QDII）"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "18.44"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "85.51"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "1.83"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.3375"
$q4.Range("H2").Value = 10

$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "014002"
$q4.Range("C3").Value = "浦银安盛全球智能科技股票（QDII）C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.30"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "42.55"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "1.24"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0037"
$q4.Range("H3").Value = 9

$q4.Range("A4").Value = 2
$q4.Range("B4").NumberFormat = "@"
$q4.Range("B4").Value = "006555"
$q4.Range("C4").Value = "浦银安盛全球智能科技股票（QDII）A"
$q4.Range("D4").NumberFormat = "@"
$q4.Range("D4").Value = "0.25"
$q4.Range("E4").NumberFormat = "@"
$q4.Range("E4").Value = "42.55"
$q4.Range("F4").NumberFormat = "@"
$q4.Range("F4").Value = "1.24"
$q4.Range("G4").NumberFormat = "@"
$q4.Range("G4").Value = "0.0031"
$q4.Range("H4").Value = 9

# Clear the now-unused rows 5:6 left over from the old 6-row sheet.
$q4.Range("A5:H6").Clear()

# ---------------------------------------------------------------------
# Step 4: update the "总计" overview sheet with a new 2022-Q4 row.
# ---------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")
$zj.Range("A2:D2").Copy($zj.Range("A3"))

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 3
$zj.Range("D2").Value = 0.34

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 5
$zj.Range("D3").Value = 1.28
